# Fixing network data cleaning scripts
# - Rename header columns to short snake_case names
# - Title-case the Spanish connector words (de/del/la/las/el/los/y) in
#   state/municipality names so they read "De"/"Del"/"La"/"Las"/"El"/"Los"/"Y"
# - Correct two floating point rounding values
# - Remove the trailing metadata/footnote rows and shrink the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Title-case the Spanish connector words inside the state (col A) and
#    municipality (col B) name cells. Walk the data rows only (header row 1
#    is handled separately below); going up to 969 covers the footer rows
#    too (harmless, they get deleted afterwards).
for ($r = 2; $r -le 969; $r++) {
    for ($col = 1; $col -le 2; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value2
        if ($v -is [string]) {
            $new = $v -replace '\bde\b', 'De' `
                       -replace '\bdel\b', 'Del' `
                       -replace '\bla\b', 'La' `
                       -replace '\blas\b', 'Las' `
                       -replace '\bel\b', 'El' `
                       -replace '\blos\b', 'Los' `
                       -replace '\by\b', 'Y'
            $cell.Value = $new
        }
    }
}

# 2) Rename the header row to short machine-friendly names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Correct the two floating point percentage values that were off by a
#    single ULP.
$ws.Cells.Item(226, 4).Value = 0.009652906140891352
$ws.Cells.Item(493, 4).Value = 0.009652906140891352

# 4) Drop the trailing "Tamaño de la muestra / Fuente / Elaborado por /
#    Secretaría / Marzo de 2017" footnote rows (965-969); this also shrinks
#    the sheet's used range/dimension down to A1:D963.
$ws.Rows("965:969").Delete()
